# Edit corresponding to: updating the interest-rate input on the Q2
# worksheet and clearing the formula in the first row of the discounting
# table (I9), then leaving the Q2 sheet/cell as the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Q2")

# Change the nominal discount rate input from 0.004 to 0.009.
$ws.Range("I4").Value = 0.009

# Clear the shared formula in I9, replacing it with a literal 0 (matches
# H9 = 0 contribution period, same result as before but now a hard value).
$ws.Range("I9").Value = 0

# Make Q2 the active sheet and select I10, matching the final view state.
$ws.Activate()
$ws.Range("I10").Select()
